$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.413.99"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "1.641.00"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'1.003"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'303.64"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'0.3771"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'52.31"
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("D9").Value = "'0.3657"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").Value = "'1.249"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").Value = "'0.08117"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "'22.95"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "'6.647"
$ws.Range("E14").Value = "  +0.90%  "
$ws.Range("D15").Value = "'0.00001256"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "'7.316"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "1.641.57"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "'94.33"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "'0.06948"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'18.17"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'6.560"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "23.425.96"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "'12.87"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'3.245"
$ws.Range("D26").Value = "'2.446"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "'21.26"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "'151.52"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "'5.313"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "'136.32"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "'2.322"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").Value = "1.824.74"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").Value = "'6.912"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").Value = "'10.99"
$ws.Range("E34").Value = "  +7.12%  "
$ws.Range("D35").Value = "'0.9641"
$ws.Range("E35").Value = "  +0.91%  "
$ws.Range("D36").Value = "'0.02864"
$ws.Range("E36").Value = "  +3.64%  "
$ws.Range("D37").Value = "'6.288"
$ws.Range("E37").Value = "  +3.08%  "
$ws.Range("D38").Value = "'0.2563"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").Value = "'0.07304"
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("D40").Value = "'0.08853"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'1.377"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("D42").Value = "'0.7138"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'16.44"
$ws.Range("E43").Value = "  +4.54%  "
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("D45").Value = "'0.6567"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "'2.366"
$ws.Range("E46").Value = "  +1.80%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "'3.999"
$ws.Range("D49").Value = "'0.08002"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "'1.218"
$ws.Range("E50").Value = "  +1.68%  "
$ws.Range("D51").Value = "'128.36"
$ws.Range("E51").Value = "  -4.11%  "
